# Weekly update: insert a new price record for "Acelga" at Macroferia
# Regional de Talca, shifting the existing history rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 106 (rows 106:164 shift down to 107:165).
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row with the latest weekly observation.
$ws.Range("A106").Value = 5
$ws.Range("B106").Value = "Macroferia Regional de Talca"
$ws.Range("C106").Value = "Maule"
$ws.Range("D106").Value = 44452
$ws.Range("E106").Value = 7
$ws.Range("F106").Value = 100112009
$ws.Range("G106").Value = "Acelga"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 400
$ws.Range("K106").Value = 2500
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = 2500
$ws.Range("N106").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O106").Value = "Región del Maule"
$ws.Range("P106").Value = 625
$ws.Range("Q106").Value = 4
$ws.Range("R106").Value = "Hortaliza"
